$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled like the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column
$values = @(1, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
